$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the L-series sample IDs (rows 2-25) to the LISO-prefixed scheme
$ws.Range("A2").Value = "LISO1"
$ws.Range("A3").Value = "LISO2"
$ws.Range("A4").Value = "LISO3"
$ws.Range("A5").Value = "LISO4"
$ws.Range("A6").Value = "LISO5"
$ws.Range("A7").Value = "LISO6"
$ws.Range("A8").Value = "LISO7"
$ws.Range("A9").Value = "LISO8"
$ws.Range("A10").Value = "LISO9"
$ws.Range("A11").Value = "LISO10"
$ws.Range("A12").Value = "LISO11"
$ws.Range("A13").Value = "LISO12"
$ws.Range("A14").Value = "LISO13"
$ws.Range("A15").Value = "LISO14"
$ws.Range("A16").Value = "LISO15"
$ws.Range("A17").Value = "LISO16"
$ws.Range("A18").Value = "LISO17"
$ws.Range("A19").Value = "LISO18"
$ws.Range("A20").Value = "LISO19"
$ws.Range("A21").Value = "LISO20"
$ws.Range("A22").Value = "LISO21"
$ws.Range("A23").Value = "LISO22"
$ws.Range("A24").Value = "LISO23"
$ws.Range("A25").Value = "LISO24"

# Rename the R-series sample IDs (rows 26-46) to the RISO-prefixed scheme
$ws.Range("A26").Value = "RISO1"
$ws.Range("A27").Value = "RISO2"
$ws.Range("A28").Value = "RISO3"
$ws.Range("A29").Value = "RISO4"
$ws.Range("A30").Value = "RISO5"
$ws.Range("A31").Value = "RISO6"
$ws.Range("A32").Value = "RISO7"
$ws.Range("A33").Value = "RISO8"
$ws.Range("A34").Value = "RISO9"
$ws.Range("A35").Value = "RISO10"
$ws.Range("A36").Value = "RISO11"
$ws.Range("A37").Value = "RISO12"
$ws.Range("A38").Value = "RISO13"
$ws.Range("A39").Value = "RISO14"
$ws.Range("A40").Value = "RISO15"
$ws.Range("A41").Value = "RISO16"
$ws.Range("A42").Value = "RISO17"
$ws.Range("A43").Value = "RISO18"
$ws.Range("A44").Value = "RISO19"
$ws.Range("A45").Value = "RISO20"
$ws.Range("A46").Value = "RISO21"

# Rename the C-series sample IDs (rows 47-74) to the CISO-prefixed scheme
$ws.Range("A47").Value = "CISO1"
$ws.Range("A48").Value = "CISO2"
$ws.Range("A49").Value = "CISO3"
$ws.Range("A50").Value = "CISO4"
$ws.Range("A51").Value = "CISO5"
$ws.Range("A52").Value = "CISO6"
$ws.Range("A53").Value = "CISO7"
$ws.Range("A54").Value = "CISO8"
$ws.Range("A55").Value = "CISO9"
$ws.Range("A56").Value = "CISO10"
$ws.Range("A57").Value = "CISO11"
$ws.Range("A58").Value = "CISO12"
$ws.Range("A59").Value = "CISO13"
$ws.Range("A60").Value = "CISO14"
$ws.Range("A61").Value = "CISO15"
$ws.Range("A62").Value = "CISO16"
$ws.Range("A63").Value = "CISO17"
$ws.Range("A64").Value = "CISO18"
$ws.Range("A65").Value = "CISO19"
$ws.Range("A66").Value = "CISO20"
$ws.Range("A67").Value = "CISO21"
$ws.Range("A68").Value = "CISO22"
$ws.Range("A69").Value = "CISO23"
$ws.Range("A70").Value = "CISO24"
$ws.Range("A71").Value = "CISO25"
$ws.Range("A72").Value = "CISO26"
$ws.Range("A73").Value = "CISO27"
$ws.Range("A74").Value = "CISO28"

# Move the active selection to E7 (matches the saved cursor position)
$ws.Range("E7").Select() | Out-Null
